$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at sheet row 157 (pushing the existing rows 157:273 down to 158:274).
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(157, 1).Value = 3
$ws.Cells.Item(157, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(157, 3).Value = "Coquimbo"
$ws.Cells.Item(157, 4).Value = 44574
$ws.Cells.Item(157, 5).Value = 5
$ws.Cells.Item(157, 6).Value = 100112043
$ws.Cells.Item(157, 7).Value = "Pepino ensalada"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 60
$ws.Cells.Item(157, 11).Value = 10000
$ws.Cells.Item(157, 12).Value = 10000
$ws.Cells.Item(157, 13).Value = 10000
$ws.Cells.Item(157, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(157, 15).Value = "Limache"
$ws.Cells.Item(157, 16).Value = 143
$ws.Cells.Item(157, 17).Value = 70
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Match the style (date number format) used by the other rows' date column.
$ws.Cells.Item(157, 4).NumberFormat = $ws.Cells.Item(158, 4).NumberFormat
